$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H96").Value = 1954.4
$ws.Range("I96").Value = 2689.1667
$ws.Range("J96").Value = 852.25
$ws.Range("K96").Value = 8067.500100000001
$ws.Range("L96").Value = 2556.75
$ws.Range("M96").Value = -6694.500100000001

$ws.Range("H97").Value = 342.8
$ws.Range("I97").Value = 500
$ws.Range("J97").Value = 303.5
$ws.Range("K97").Value = 1500
$ws.Range("L97").Value = 910.5
$ws.Range("M97").Value = -1004
$ws.Range("N97").Value = -1902.5

$ws.Range("H112").Value = 1956.421
$ws.Range("I112").Value = 1084.6666
$ws.Range("J112").Value = 2119.875
$ws.Range("K112").Value = 3253.9998
$ws.Range("L112").Value = 6359.625
$ws.Range("M112").Value = -2145.9998

$ws.Range("H132").Value = 8138378
$ws.Range("I132").Value = 9809494
$ws.Range("J132").Value = 21529.428
$ws.Range("K132").Value = 29428482
$ws.Range("L132").Value = 64588.284
$ws.Range("M132").Value = -29425952

$ws.Range("H138").Value = 2473.3445
$ws.Range("I138").Value = 1352
$ws.Range("J138").Value = 2734.4795
$ws.Range("K138").Value = 4056
$ws.Range("L138").Value = 8203.4385
$ws.Range("M138").Value = 1084
$ws.Range("N138").Value = -18483.4385

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6477.32
$ws.Range("I32").Value = 5166.106
$ws.Range("J32").Value = 13907.533
$ws.Range("K32").Value = 5166.106
$ws.Range("L32").Value = 13907.533
$ws.Range("M32").Value = -4879.106
$ws.Range("N32").Value = -14481.533

$ws.Range("H61").Value = 1507.1765
$ws.Range("I61").Value = 1009
$ws.Range("J61").Value = 2702.8
$ws.Range("K61").Value = 1009
$ws.Range("L61").Value = 2702.8
$ws.Range("M61").Value = -797
$ws.Range("N61").Value = -3126.8

$ws.Range("H74").Value = 1852.125
$ws.Range("I74").Value = 967.6957
$ws.Range("J74").Value = 4112.3335
$ws.Range("K74").Value = 967.6957
$ws.Range("L74").Value = 4112.3335
$ws.Range("M74").Value = -93.69569999999999

$ws.Range("H77").Value = 1852.125
$ws.Range("I77").Value = 967.6957
$ws.Range("J77").Value = 4112.3335
$ws.Range("K77").Value = 4838.4785
$ws.Range("L77").Value = 20561.6675
$ws.Range("M77").Value = -470.4785000000002

$ws.Range("H132").Value = 3412.2593
$ws.Range("I132").Value = 3189.2222
$ws.Range("J132").Value = 3858.3333
$ws.Range("K132").Value = 9567.6666
$ws.Range("L132").Value = 11574.9999
$ws.Range("M132").Value = -7037.6666
$ws.Range("N132").Value = -16634.9999

$ws.Range("H136").Value = 1507.1765
$ws.Range("I136").Value = 1009
$ws.Range("J136").Value = 2702.8
$ws.Range("K136").Value = 3027
$ws.Range("L136").Value = 8108.400000000001
$ws.Range("M136").Value = -477
$ws.Range("N136").Value = -13208.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4324.875
$ws.Range("I20").Value = 4000
$ws.Range("J20").Value = 4866.3335
$ws.Range("K20").Value = 4000
$ws.Range("L20").Value = 4866.3335
$ws.Range("M20").Value = -3753
$ws.Range("N20").Value = -5360.3335

$ws.Range("H94").Value = 8621494
$ws.Range("I94").Value = 12500865
$ws.Range("J94").Value = 669.8889
$ws.Range("K94").Value = 12500865
$ws.Range("L94").Value = 669.8889
$ws.Range("M94").Value = -12500414
$ws.Range("N94").Value = -1571.8889

$ws.Range("H105").Value = 126239560
$ws.Range("I105").Value = 144273500
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 144273500
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -144271753

$ws.Range("H107").Value = 1206.9
$ws.Range("I107").Value = 1141
$ws.Range("J107").Value = 1800
$ws.Range("K107").Value = 1141
$ws.Range("L107").Value = 1800
$ws.Range("M107").Value = 779

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1378.8813
$ws.Range("I31").Value = 1237.2452
$ws.Range("J31").Value = 2630
$ws.Range("K31").Value = 1237.2452
$ws.Range("L31").Value = 2630
$ws.Range("M31").Value = -942.2452000000001
$ws.Range("N31").Value = -3220

$ws.Range("H34").Value = 1378.8813
$ws.Range("I34").Value = 1237.2452
$ws.Range("J34").Value = 2630
$ws.Range("K34").Value = 1237.2452
$ws.Range("L34").Value = 2630
$ws.Range("M34").Value = -1035.2452
$ws.Range("N34").Value = -3034

$ws.Range("H58").Value = 8843.200000000001
$ws.Range("I58").Value = 1358.5
$ws.Range("J58").Value = 13833
$ws.Range("K58").Value = 1358.5
$ws.Range("L58").Value = 13833
$ws.Range("M58").Value = -1155.5
$ws.Range("N58").Value = -14239

$ws.Range("H86").Value = 3055424.2
$ws.Range("I86").Value = 5146002.5
$ws.Range("J86").Value = 35700.332
$ws.Range("K86").Value = 5146002.5
$ws.Range("L86").Value = 35700.332
$ws.Range("M86").Value = -5144879.5
$ws.Range("N86").Value = -37946.332

$ws.Range("H89").Value = 3055424.2
$ws.Range("I89").Value = 5146002.5
$ws.Range("J89").Value = 35700.332
$ws.Range("K89").Value = 25730012.5
$ws.Range("L89").Value = 178501.66
$ws.Range("M89").Value = -25724396.5
$ws.Range("N89").Value = -189733.66

$ws.Range("H132").Value = 1771.0344
$ws.Range("I132").Value = 1491.75
$ws.Range("J132").Value = 2391.6667
$ws.Range("K132").Value = 4475.25
$ws.Range("L132").Value = 7175.000100000001
$ws.Range("M132").Value = -1945.25
$ws.Range("N132").Value = -12235.0001

$ws.Range("H134").Value = 14707294
$ws.Range("I134").Value = 1429.3529
$ws.Range("J134").Value = 29413158
$ws.Range("K134").Value = 4288.0587
$ws.Range("L134").Value = 88239474
$ws.Range("M134").Value = -1753.0587
$ws.Range("N134").Value = -88244544

$ws.Range("H136").Value = 8843.200000000001
$ws.Range("I136").Value = 1358.5
$ws.Range("J136").Value = 13833
$ws.Range("K136").Value = 4075.5
$ws.Range("L136").Value = 41499
$ws.Range("M136").Value = -1525.5
$ws.Range("N136").Value = -46599

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 463.6316
$ws.Range("I5").Value = 378
$ws.Range("J5").Value = 2005
$ws.Range("K5").Value = 1134
$ws.Range("L5").Value = 6015
$ws.Range("M5").Value = -1022
$ws.Range("N5").Value = -6239

$ws.Range("H58").Value = 1541.75
$ws.Range("I58").Value = 497.5
$ws.Range("J58").Value = 1750.6
$ws.Range("K58").Value = 1492.5
$ws.Range("L58").Value = 5251.799999999999
$ws.Range("M58").Value = -1364.5
$ws.Range("N58").Value = -5507.799999999999

$ws.Range("H81").Value = 2624.75
$ws.Range("I81").Value = 2624.75
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 7874.25
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -6751.25
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 2624.75
$ws.Range("I84").Value = 2624.75
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 23622.75
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -18006.75
$ws.Range("N84").ClearContents()

$ws.Range("H87").Value = 5444
$ws.Range("I87").Value = 888
$ws.Range("J87").Value = 10000
$ws.Range("K87").Value = 2664
$ws.Range("L87").Value = 30000
$ws.Range("M87").Value = -1416
$ws.Range("N87").Value = -32496

$ws.Range("H90").Value = 5444
$ws.Range("I90").Value = 888
$ws.Range("J90").Value = 10000
$ws.Range("K90").Value = 7992
$ws.Range("L90").Value = 90000
$ws.Range("M90").Value = -1752
$ws.Range("N90").Value = -102480

$ws.Range("H97").Value = 1100
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 1100
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 3300
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -4292

$ws.Range("H131").Value = 32307378
$ws.Range("I131").Value = 66667036
$ws.Range("J131").Value = 95196
$ws.Range("K131").Value = 200001108
$ws.Range("L131").Value = 285588
$ws.Range("M131").Value = -199996068
$ws.Range("N131").Value = -295668

$ws.Range("H135").Value = 463.6316
$ws.Range("I135").Value = 378
$ws.Range("J135").Value = 2005
$ws.Range("K135").Value = 3402
$ws.Range("L135").Value = 18045
$ws.Range("M135").Value = -867
$ws.Range("N135").Value = -23115

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 50003040
$ws.Range("I70").Value = 35717904
$ws.Range("J70").Value = 100001010
$ws.Range("K70").Value = 35717904
$ws.Range("L70").Value = 100001010
$ws.Range("M70").Value = -35717634
$ws.Range("N70").Value = -100001550

$ws.Range("H73").Value = 50003040
$ws.Range("I73").Value = 35717904
$ws.Range("J73").Value = 100001010
$ws.Range("K73").Value = 35717904
$ws.Range("L73").Value = 100001010
$ws.Range("M73").Value = -35716968
$ws.Range("N73").Value = -100002882

$ws.Range("H97").Value = 1127.9333
$ws.Range("I97").Value = 1028.909
$ws.Range("J97").Value = 1400.25
$ws.Range("K97").Value = 1028.909
$ws.Range("L97").Value = 1400.25
$ws.Range("M97").Value = -532.9090000000001

$ws.Range("H132").Value = 16524.334
$ws.Range("I132").Value = 31155.5
$ws.Range("J132").Value = 4819.4
$ws.Range("K132").Value = 93466.5
$ws.Range("L132").Value = 14458.2
$ws.Range("M132").Value = -90936.5
$ws.Range("N132").Value = -19518.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 939
$ws.Range("I22").Value = 922.4
$ws.Range("J22").Value = 966.6667
$ws.Range("K22").Value = 922.4
$ws.Range("L22").Value = 966.6667
$ws.Range("M22").Value = -627.4
$ws.Range("N22").Value = -1556.6667

$ws.Range("H27").Value = 939
$ws.Range("I27").Value = 922.4
$ws.Range("J27").Value = 966.6667
$ws.Range("K27").Value = 922.4
$ws.Range("L27").Value = 966.6667
$ws.Range("M27").Value = -815.4
$ws.Range("N27").Value = -1180.6667

$ws.Range("H46").Value = 3738
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 3738
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 3738
$ws.Range("N46").Value = -4114
